# Projet 3.pptx edits
#  1. Slide 2 ("Sommaire"): split the "Evènement gardien" bullet into two
#     runs: "Evènements " + "gardien".
#  2. Slide 3 ("Présentation du projet"): drop the automatic line-space
#     reduction on the content placeholder (normAutofit lnSpcReduction
#     -> normAutofit).
#  3. Slide 3: remove the "Librairie utilisée" bullet paragraph
#     (lvl 1) from the content placeholder.

$p = $ppt.ActivePresentation

# --- 1. Slide 2: "Evènement gardien" -> "Evènements " / "gardien" ---
$slide2 = $p.Slides.Item(2)
$contentShape2 = $slide2.Shapes.Item(2)
$tr2 = $contentShape2.TextFrame.TextRange
for ($i = 1; $i -le $tr2.Count; $i++) {
    $candidate = $tr2.Paragraphs($i)
    if ($candidate.Text.TrimEnd("`r") -eq "Evènement gardien") {
        # Rewrite only the first 10 characters ("Evènement ") so the engine
        # keeps the change minimal and splits the run right after the
        # inserted "s", leaving "gardien" as its own run.
        $candidate.Characters(1, 10).Text = "Evènements "
        break
    }
}

# --- 2/3. Slide 3: content placeholder tweaks ---
$slide3 = $p.Slides.Item(3)
$contentShape3 = $slide3.Shapes.Item(2)

# Reset AutoSize so the normAutofit element no longer carries the
# lnSpcReduction="10000" attribute (regenerates a bare <a:normAutofit/>).
$contentShape3.TextFrame.AutoSize = 2

# Remove the "Librairie utilisée" paragraph entirely.
# (TextRange.Text includes the trailing paragraph mark "`r", so trim it
# before comparing against the literal bullet text.)
$tr3 = $contentShape3.TextFrame.TextRange
for ($i = 1; $i -le $tr3.Count; $i++) {
    $para = $tr3.Paragraphs($i)
    if ($para.Text.TrimEnd("`r") -eq "Librairie utilisée") {
        $para.Delete(1)
        break
    }
}
